$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style of the existing
# header cells (e.g. G1 / "sum").
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding data value in H2.
$ws.Range("H2").Value = 1
